$d = $word.ActiveDocument

# The "Conclusion" section ends with a paragraph containing the sentence
# "The Board should affirm the MAC's adjustment."  It is followed by two
# empty spacer paragraphs right before the section properties. Those two
# empty paragraphs are being removed so the section properties follow
# directly after the "Conclusion" paragraph.

$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*The Board should affirm the MAC*adjustment*") {
        $targetIdx = $i
    }
}

if ($targetIdx -eq -1) {
    throw "Could not locate the 'Board should affirm' paragraph"
}

$count = $d.Paragraphs.Count

if ($targetIdx -le $count - 2) {
    $first = $d.Paragraphs.Item($targetIdx + 1)
    $second = $d.Paragraphs.Item($targetIdx + 2)

    # Step 1: remove the first trailing empty paragraph. Deleting the range
    # from its start to the start of the next paragraph merges it forward
    # and cleanly drops its own paragraph mark/formatting, leaving the
    # second paragraph's mark intact.
    $r1 = $d.Range($first.Range.Start, $second.Range.Start)
    $r1.Delete()

    # Step 2: the former "second" empty paragraph is now the very last
    # paragraph in the document. Remove its mark as well by extending the
    # deleted range one character back (into the end of the target
    # paragraph's own mark) so no stray empty paragraph is left behind.
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r2 = $d.Range($last.Range.Start - 1, $last.Range.End)
    $r2.Delete()
}
